# Update countries & provincias Spain
# - Refresh the "last updated" timestamp string.
# - Refresh COVID-19 per-country counters for the rows that changed
#   (India, Belgica, Indonesia, Suiza, Oman, Rumania, Moldavia, Marruecos,
#   Malasia, Finlandia, Albania, Eslovenia, Hong Kong, Namibia).
#   Suiza/Oman also swap rank order (row 41 now shows Oman's bigger count,
#   row 42 now shows Suiza's old count) while each row keeps its own
#   country name cell updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Junio de 2020 a las 12:08"

# Row 7 - India
$ws.Range("B7").Value = 441643
$ws.Range("C7").Value = 1193
$ws.Range("D7").Value = 248450
$ws.Range("E7").Value = 179166
$ws.Range("G7").Value = 12
$ws.Range("H7").Value = 14027

# Row 26 - Belgica
$ws.Range("B26").Value = 60810
$ws.Range("C26").Value = 260
$ws.Range("E26").Value = 34326
$ws.Range("G26").Value = 17
$ws.Range("H26").Value = 9713

# Row 32 - Indonesia
$ws.Range("B32").Value = 47896
$ws.Range("C32").Value = 1051
$ws.Range("D32").Value = 19241
$ws.Range("E32").Value = 26120
$ws.Range("G32").Value = 35
$ws.Range("H32").Value = 2535

# Row 41 - now Oman (name + data)
$ws.Range("A41").Value = "Oman"
$ws.Range("B41").Value = 32394
$ws.Range("C41").Value = 1318
$ws.Range("D41").Value = 17279
$ws.Range("E41").Value = 14975
$ws.Range("G41").Value = 3
$ws.Range("H41").Value = 140

# Row 42 - now Suiza (name + data)
$ws.Range("A42").Value = "Suiza"
$ws.Range("B42").Value = 31310
$ws.Range("D42").Value = 29000
$ws.Range("E42").Value = 354
$ws.Range("H42").Value = 1956

# Row 49 - Rumania
$ws.Range("B49").Value = 24505
$ws.Range("C49").Value = 214
$ws.Range("D49").Value = 17187
$ws.Range("E49").Value = 5779
$ws.Range("G49").Value = 16
$ws.Range("H49").Value = 1539

# Row 57 - Moldavia
$ws.Range("E57").Value = 5860
$ws.Range("G57").Value = 4
$ws.Range("H57").Value = 484

# Row 68 - Marruecos
$ws.Range("B68").Value = 10264
$ws.Range("C68").Value = 92
$ws.Range("D68").Value = 8384
$ws.Range("E68").Value = 1666

# Row 72 - Malasia
$ws.Range("B72").Value = 8590
$ws.Range("C72").Value = 3
$ws.Range("D72").Value = 8186
$ws.Range("E72").Value = 283

# Row 75 - Finlandia
$ws.Range("B75").Value = 7155
$ws.Range("C75").Value = 11
$ws.Range("E75").Value = 428

# Row 105 - Albania
$ws.Range("B105").Value = 2047
$ws.Range("C105").Value = 52
$ws.Range("D105").Value = 1195
$ws.Range("E105").Value = 807
$ws.Range("G105").Value = 1
$ws.Range("H105").Value = 45

# Row 118 - Eslovenia
$ws.Range("B118").Value = 1534
$ws.Range("C118").Value = 13
$ws.Range("E118").Value = 49

# Row 123 - Hong Kong
$ws.Range("B123").Value = 1178
$ws.Range("C123").Value = 16
$ws.Range("D123").Value = 1083
$ws.Range("E123").Value = 89
$ws.Range("G123").Value = 1
$ws.Range("H123").Value = 6

# Row 186 - Namibia
$ws.Range("B186").Value = 67
$ws.Range("C186").Value = 4
$ws.Range("E186").Value = 46
